# taskkill support +1C coloring for &
# Adds the translator "Dieter Hummel" (German) with mailto hyperlink to the
# "ааа" worksheet, right below the existing entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Values are entered Name, Language, email - matching how the new shared
# strings were appended (Dieter Hummel, German, dh@level47.de) in the
# underlying workbook XML.
$ws.Cells.Item($row, 2).Value = "Dieter Hummel"
$ws.Cells.Item($row, 1).Value = "German"
$ws.Cells.Item($row, 3).Value = "dh@level47.de "

# Add the mailto hyperlink on the email cell, mirroring the other rows in
# the sheet (which use the "hyperlink" cell style).
$emailCell = $ws.Cells.Item($row, 3)
$ws.Hyperlinks.Add($emailCell, "mailto:dh@level47.de") | Out-Null
$emailCell.Style = "Гиперссылка"

# Move the active selection down to the next empty row, as happens after
# entering a new row of data in Excel.
$ws.Range("C7").Select()
